$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) and Volume(1h) (E) updates ---
$ws.Range("D2").Value = '68.531.67'
$ws.Range("E2").Value = '  +0.51%  '

$ws.Range("D3").Value = '2.694.45'
$ws.Range("E3").Value = '  +1.93%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '''598.52'
$ws.Range("E5").Value = '  +0.18%  '

$ws.Range("D6").Value = '''159.45'
$ws.Range("E6").Value = '  +2.12%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '''0.543'
$ws.Range("E8").Value = '  -0.27%  '

$ws.Range("D9").Value = '2.692.84'
$ws.Range("E9").Value = '  +1.85%  '

$ws.Range("E10").Value = '  -3.86%  '

$ws.Range("E11").Value = '  -0.85%  '

$ws.Range("E12").Value = '  +0.86%  '

$ws.Range("E13").Value = '  +1.72%  '

$ws.Range("D14").Value = '''28.31'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '3.209.30'
$ws.Range("E15").Value = '  +2.61%  '

$ws.Range("D16").Value = '''0.0000187'
$ws.Range("E16").Value = '  -2.13%  '

$ws.Range("D17").Value = '68.455.11'
$ws.Range("E17").Value = '  +0.23%  '

$ws.Range("D18").Value = '2.686.78'
$ws.Range("E18").Value = '  +1.42%  '

$ws.Range("D19").Value = '''11.88'
$ws.Range("E19").Value = '  +4.39%  '

$ws.Range("D20").Value = '''366.14'
$ws.Range("E20").Value = '  +0.79%  '

$ws.Range("D21").Value = '''7.58'
$ws.Range("E21").Value = '  +1.79%  '

$ws.Range("D22").Value = '''4.53'
$ws.Range("E22").Value = '  +2.95%  '

$ws.Range("D23").Value = '''4.93'
$ws.Range("E23").Value = '  +2.07%  '

$ws.Range("D24").Value = '''2.13'
$ws.Range("E24").Value = '  +3.39%  '

$ws.Range("D25").Value = '''74.98'
$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = '''0.998'
$ws.Range("E26").Value = '  -0.18%  '

$ws.Range("D27").Value = '''10.16'
$ws.Range("E27").Value = '  +4.20%  '

$ws.Range("E29").Value = '  -0.98%  '

$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").Value = '''579.29'
$ws.Range("E31").Value = '  +4.14%  '

$ws.Range("D32").Value = '''8.25'
$ws.Range("E32").Value = '  +3.00%  '

$ws.Range("E33").Value = '  +1.41%  '

$ws.Range("E34").Value = '  +4.30%  '

$ws.Range("D35").Value = '''1.65'
$ws.Range("E35").Value = '  +5.44%  '

$ws.Range("E36").Value = '  +1.46%  '

$ws.Range("E37").Value = '  -0.03%  '

$ws.Range("E40").Value = '  +1.95%  '

$ws.Range("D41").Value = '''1.91'
$ws.Range("E41").Value = '  +1.81%  '

$ws.Range("D42").Value = '''5.42'
$ws.Range("E42").Value = '  +1.65%  '

$ws.Range("D43").Value = '''17.87'
$ws.Range("E43").Value = '  +0.53%  '

$ws.Range("D44").Value = '''2.65'
$ws.Range("E44").Value = '  +0.44%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("E46").Value = '  -7.25%  '

$ws.Range("D47").Value = '''158.52'
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").Value = '''3.92'
$ws.Range("E48").Value = '  +5.04%  '

$ws.Range("D49").Value = '''1.77'
$ws.Range("E49").Value = '  +4.28%  '

$ws.Range("E50").Value = '  +7.06%  '

$ws.Range("D51").Value = '''22.09'
$ws.Range("E51").Value = '  +0.51%  '

# --- Rows 38/39 swapped (Monero <-> EthereumClassic) with updated values ---
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = '''20.06'
$ws.Range("E38").Value = '  +3.33%  '

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = '''161.76'
$ws.Range("E39").Value = '  +0.83%  '
